# LEonard Design Rev 13 - manual UI cleanups (small)
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$shape = $s.Shapes.Item(2)
$tf = $shape.TextFrame
$tr = $tf.TextRange

# NOTE: edits are applied from the end of the text backwards so that the
# character offsets (taken from the original, unmodified text) for the
# earlier edits stay valid.

# 1) Remove the now-redundant "<tab>"-only paragraph (old paragraph 12),
#    merging it away so only a single trailing empty paragraph remains.
$tabPara = $tr.Characters(499, 2)
$tabPara.Delete()

# 2) Collapse the "LEonardRoot" bullet's trailing runs ("Flder"/" where "/
#    "LEonard is running ") into one cleaned-up description.
$rootSub = $tr.Characters(465, 33)
$rootSub.Text = ": Folder for the LEonard directory " + [char]8220 + "tree" + [char]8221 + " " + [char]9

# 3) "All LEVariables ... " bullet: " duplicated" -> " are duplicated"
$dupSub = $tr.Characters(116, 30)
$dupSub.Text = " are duplicated to Java and Python"

# 4) Let PowerPoint recompute the body placeholder's shrink-to-fit state now
#    that the text is shorter (drops the stale lnSpcReduction="10000").
$tf.AutoSize = 2
